$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 253 (shifts rows 253..323 down to 254..324)
$ws.Range("A253:R253").Insert(-4121)

# The newly inserted blank row 253 should receive the data that row 252
# currently holds (the row data cascades down by one position).
$ws.Range("A252:R252").Copy($ws.Range("A253:R253"))

# Now overwrite row 252 with the new/updated record values.
$ws.Range("D252").Value = 44924
$ws.Range("J252").Value = 90
$ws.Range("K252").Value = 14000
$ws.Range("L252").Value = 16000
$ws.Range("M252").Value = 15000
$ws.Range("O252").Value = "Región de Arica y Parinacota"
$ws.Range("P252").Value = 300
